$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Backlog")

# Remove the "now"/"noow" status markers that had been placed in column E
$ws.Range("E28").ClearContents()
$ws.Range("E29").ClearContents()
$ws.Range("E30").ClearContents()
$ws.Range("E31").ClearContents()
$ws.Range("E32").ClearContents()
$ws.Range("E33").ClearContents()
$ws.Range("E41").ClearContents()

# Remove the "done" marker that had been placed in column F
$ws.Range("F38").ClearContents()

# Match the saved selection / scroll position
$ws.Range("F38").Select()
